# Applies the crypto price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.370.98"
$ws.Range("E2").Value = "  +0.24%  "
# Row 3
$ws.Range("D3").Value = "2.237.70"
$ws.Range("E3").Value = "  -0.44%  "
# Row 4
$ws.Range("E4").Value = "  -0.06%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.42"
$ws.Range("E5").Value = "  -0.75%  "
# Row 6
$ws.Range("E6").Value = "  -0.67%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.53"
$ws.Range("E7").Value = "  -3.21%  "
# Row 8
$ws.Range("E8").Value = "  -0.01%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.622"
$ws.Range("E9").Value = "  -0.70%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.16"
$ws.Range("E10").Value = "  +2.36%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0965"
$ws.Range("E11").Value = "  +1.27%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.14"
$ws.Range("E12").Value = "  +0.36%  "
# Row 13
$ws.Range("E13").Value = "  -0.18%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.49"
$ws.Range("E14").Value = "  -1.82%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.853"
$ws.Range("E15").Value = "  -0.79%  "
# Row 16
$ws.Range("D16").Value = "2.216.62"
$ws.Range("E16").Value = "  -1.26%  "
# Row 17
$ws.Range("D17").Value = "42.260.25"
$ws.Range("E17").Value = "  +0.41%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000111"
$ws.Range("E18").Value = "  +12.33%  "
# Row 19
$ws.Range("E19").Value = "  +1.14%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.16"
$ws.Range("E20").Value = "  +0.27%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.28"
$ws.Range("E21").Value = "  +36.12%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "231.40"
$ws.Range("E22").Value = "  -0.20%  "
# Row 23
$ws.Range("E23").Value = "  -4.58%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.74"
$ws.Range("E24").Value = "  +3.65%  "
# Row 25
$ws.Range("E25").Value = "  +0.03%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.65"
$ws.Range("E26").Value = "  +1.20%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.29"
$ws.Range("E27").Value = "  +0.06%  "
# Row 28
$ws.Range("E28").Value = "  +5.55%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.68"
$ws.Range("E29").Value = "  -2.01%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.95"
$ws.Range("E30").Value = "  +1.83%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.87"
$ws.Range("E31").Value = "  +19.00%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0806"
$ws.Range("E32").Value = "  -3.65%  "
# Row 33
$ws.Range("E33").Value = "  -1.75%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "29.77"
$ws.Range("E34").Value = "  -8.84%  "
# Row 35
$ws.Range("E35").Value = "  -0.52%  "
# Row 36
$ws.Range("E36").Value = "  -0.75%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0309"
$ws.Range("E37").Value = "  +2.31%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.26"
$ws.Range("E38").Value = "  -8.02%  "
# Row 39
$ws.Range("E39").Value = "  -0.88%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.65"
$ws.Range("E40").Value = "  -4.16%  "
# Row 41
$ws.Range("E41").Value = "  +3.90%  "
# Row 42
$ws.Range("E42").Value = "  -0.63%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.83"
$ws.Range("E43").Value = "  +1.76%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "105.67"
$ws.Range("E44").Value = "  -6.16%  "
# Row 45
$ws.Range("E45").Value = "  +2.93%  "
# Row 46
$ws.Range("E46").Value = "  -0.34%  "
# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.39"
$ws.Range("E47").Value = "  +3.65%  "
# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.14"
$ws.Range("E48").Value = "  +0.21%  "
# Row 49
$ws.Range("E49").Value = "  +0.87%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.72"
$ws.Range("E50").Value = "  +0.92%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.08"
$ws.Range("E51").Value = "  -2.45%  "
